# Applies "All get routes added" changes to NouvellesAdressesAPI workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 1. Fix existing speed-loss route (row 9): "speedloss" -> "speedlosses"
$ws.Range("B9").Formula = '=_xlfn.CONCAT($B$2,"speedlosses/{po}/{productionLine}")'

# 2. New row 23: getSpeedLosses route (its shared string is inserted first)
$ws.Range("A23").Value = "Route::get('getSpeedLosses/{site}/{productionLine}/{beginningDate}/{endingDate}', FormController::class.'@getSpeedLosses');"
$ws.Range("B23").Formula = '=_xlfn.CONCAT($B$2,"getSpeedLosses/{site}/{productionLine}/{startingDate}/{endingDate}")'

# 3. New row 22: performance route
$ws.Range("A22").Value = "Route::get('performance/{PO}', FormController::class.'@getPerformanceForASite');"
$ws.Range("B22").Formula = '=_xlfn.CONCAT($B$2,"performance/{PO}")'

$wb.Application.Calculate()
